$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inputs")
$ws.Activate()

$ws.Range("D3").Value = 0.15
$ws.Range("E3").Value = 0.081
$ws.Range("G3").Value = 22359

$ws.Range("H4").Value = 0.23
$ws.Range("H5").Value = 0.77

$ws.Range("E6").Value = 0.06
$ws.Range("E7").Value = 0.06

$ws.Range("D8").Value = 0.68
$ws.Range("E8").Value = 0.14
$ws.Range("H8").Value = 1

$ws.Range("D9").Value = 0.035
$ws.Range("E9").Value = 0.023
$ws.Range("H9").Value = 0.069

$ws.Range("D10").Value = 0.035
$ws.Range("E10").Value = 0.023
$ws.Range("H10").Formula = "=1-H9"

$ws.Range("D11").Value = 0.24
$ws.Range("E11").Value = 0.045
$ws.Range("H11").Value = 0.93

$ws.Range("D12").Value = 0.24
$ws.Range("E12").Value = 0.045
$ws.Range("H12").Value = 0.07

$ws.Range("D13").Value = 0.55
$ws.Range("E13").Value = 0.023

$ws.Range("E14").Value = 0.1

$ws.Range("E15").Value = 0.03

$ws.Range("D16").Value = 3831
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 49138898.203365818
$ws.Range("G16").Value = 921

$ws.Range("B10").Select()
